$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper "constants"
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# Stable reference/anchor cells (row 14 is untouched by this edit) used to
# copy the correct number-format / style onto cells whose underlying type
# changes (number <-> text).
#   C14 -> text "0"      (style 14)
#   E14 -> text "***.*"  (style 14)
#   J14 -> plain integer (style 15)
#   K14 -> decimal pct   (style 16)

function Set-TextCell($addr, $text, $anchor) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($anchor).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

function Set-NumberCell($addr, $number, $anchor) {
    $ws.Range($addr).Value = $number
    $ws.Range($anchor).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Shared / rich-text strings: "Volume 30  Number 31" -> "... Number 32"
#                              "Report Covering the Week 7/31/2023 Through
#                               8/6/2023" -> "... 8/7/2023 Through 8/13/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextCell "C15" "0" "C14"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -60
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = -4.166666666666

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -19.047619047619
$ws.Range("I16").Value = 198
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = 22.981366459627
$ws.Range("L16").Value = 115.217391304348
$ws.Range("M16").Value = 52.307692307692
$ws.Range("N16").Value = -72.268907563025

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 57.142857142857
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -35.714285714285
$ws.Range("I17").Value = 281
$ws.Range("J17").Value = 204
$ws.Range("K17").Value = 37.745098039215
$ws.Range("L17").Value = 130.327868852459
$ws.Range("M17").Value = 106.617647058824
$ws.Range("N17").Value = 27.727272727272

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 41
$ws.Range("H18").Value = 32.258064516129
$ws.Range("I18").Value = 349
$ws.Range("J18").Value = 276
$ws.Range("K18").Value = 26.449275362318
$ws.Range("L18").Value = 71.078431372549
$ws.Range("M18").Value = 20.761245674740
$ws.Range("N18").Value = -76.686706746827

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 40
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 105
$ws.Range("G19").Value = 124
$ws.Range("H19").Value = -15.322580645161
$ws.Range("I19").Value = 859
$ws.Range("J19").Value = 917
$ws.Range("K19").Value = -6.324972737186
$ws.Range("L19").Value = 109.512195121951
$ws.Range("M19").Value = 118.575063613232
$ws.Range("N19").Value = 4.501216545012

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 22.222222222222
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 62.5
$ws.Range("I20").Value = 322
$ws.Range("J20").Value = 138
$ws.Range("K20").Value = 133.333333333333
$ws.Range("L20").Value = 153.543307086614
$ws.Range("M20").Value = 86.127167630057
$ws.Range("N20").Value = -86.989898989899

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold row, styles 18/19, values only)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 231
$ws.Range("G21").Value = 247
$ws.Range("H21").Value = -6.477732793522
$ws.Range("I21").Value = 2032
$ws.Range("J21").Value = 1723
$ws.Range("K21").Value = 17.933836331979
$ws.Range("L21").Value = 108.838643371017
$ws.Range("M21").Value = 78.873239436619
$ws.Range("N21").Value = -64.679297757691

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-TextCell "C22" "0" "C14"
Set-NumberCell "D22" 2 "J14"
Set-NumberCell "E22" -100 "K14"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 0

# ---------------------------------------------------------------------------
# Row 23 - Housing (only M23 changes)
# ---------------------------------------------------------------------------
$ws.Range("M23").Value = 25

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 73
$ws.Range("D24").Value = 55
$ws.Range("E24").Value = 32.727272727272
$ws.Range("F24").Value = 271
$ws.Range("G24").Value = 207
$ws.Range("H24").Value = 30.917874396135
$ws.Range("I24").Value = 1756
$ws.Range("J24").Value = 1729
$ws.Range("K24").Value = 1.561596298438
$ws.Range("L24").Value = 51.509922346850
$ws.Range("M24").Value = 91.912568306010

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -19.047619047619
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = 13.636363636363
$ws.Range("I25").Value = 533
$ws.Range("J25").Value = 429
$ws.Range("K25").Value = 24.242424242424
$ws.Range("L25").Value = 83.161512027491
$ws.Range("M25").Value = 15.618221258134

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextCell "C26" "0" "C14"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -40
$ws.Range("J26").Value = 33
$ws.Range("K26").Value = 6.060606060606

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 72
$ws.Range("J27").Value = 51
$ws.Range("K27").Value = 41.176470588235
$ws.Range("L27").Value = 75.609756097561

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic. (only L28 / N28 change)
# ---------------------------------------------------------------------------
$ws.Range("L28").Value = 20
$ws.Range("N28").Value = -64.705882352941

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc. (only L29 / N29 change)
# ---------------------------------------------------------------------------
$ws.Range("L29").Value = 50
$ws.Range("N29").Value = -57.142857142857

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-TextCell "D30" "0" "C14"
Set-TextCell "E30" "***.*" "E14"
$ws.Range("G30").Value = 1

Write-Host "Edit script complete"
